$d = $word.ActiveDocument

$d.Content.Find.Execute("83×53=", $true, $false, $false, $false, $false, $true, 1, $false, "60×87=", 2) | Out-Null
$d.Content.Find.Execute("93×90=", $true, $false, $false, $false, $false, $true, 1, $false, "60×86=", 2) | Out-Null
$d.Content.Find.Execute("56×24=", $true, $false, $false, $false, $false, $true, 1, $false, "39×75=", 2) | Out-Null
$d.Content.Find.Execute("18×62=", $true, $false, $false, $false, $false, $true, 1, $false, "84×85=", 2) | Out-Null
$d.Content.Find.Execute("97×57=", $true, $false, $false, $false, $false, $true, 1, $false, "96×49=", 2) | Out-Null
$d.Content.Find.Execute("55×74=", $true, $false, $false, $false, $false, $true, 1, $false, "15×65=", 2) | Out-Null
$d.Content.Find.Execute("27×55=", $true, $false, $false, $false, $false, $true, 1, $false, "42×96=", 2) | Out-Null
$d.Content.Find.Execute("55×68=", $true, $false, $false, $false, $false, $true, 1, $false, "85×90=", 2) | Out-Null
$d.Content.Find.Execute("18×74=", $true, $false, $false, $false, $false, $true, 1, $false, "50×86=", 2) | Out-Null
$d.Content.Find.Execute("30×48=", $true, $false, $false, $false, $false, $true, 1, $false, "27×44=", 2) | Out-Null
$d.Content.Find.Execute("92×51=", $true, $false, $false, $false, $false, $true, 1, $false, "11×70=", 2) | Out-Null
$d.Content.Find.Execute("31×66=", $true, $false, $false, $false, $false, $true, 1, $false, "65×91=", 2) | Out-Null
$d.Content.Find.Execute("85×63=", $true, $false, $false, $false, $false, $true, 1, $false, "70×41=", 2) | Out-Null
$d.Content.Find.Execute("27×40=", $true, $false, $false, $false, $false, $true, 1, $false, "17×14=", 2) | Out-Null
$d.Content.Find.Execute("26×25=", $true, $false, $false, $false, $false, $true, 1, $false, "62×70=", 2) | Out-Null
$d.Content.Find.Execute("56×68=", $true, $false, $false, $false, $false, $true, 1, $false, "12×85=", 2) | Out-Null
$d.Content.Find.Execute("61×97=", $true, $false, $false, $false, $false, $true, 1, $false, "27×88=", 2) | Out-Null
$d.Content.Find.Execute("24×58=", $true, $false, $false, $false, $false, $true, 1, $false, "22×49=", 2) | Out-Null
$d.Content.Find.Execute("47×89=", $true, $false, $false, $false, $false, $true, 1, $false, "82×93=", 2) | Out-Null
$d.Content.Find.Execute("43×19=", $true, $false, $false, $false, $false, $true, 1, $false, "79×73=", 2) | Out-Null
$d.Content.Find.Execute("59×79=", $true, $false, $false, $false, $false, $true, 1, $false, "83×30=", 2) | Out-Null
$d.Content.Find.Execute("15×57=", $true, $false, $false, $false, $false, $true, 1, $false, "44×81=", 2) | Out-Null
$d.Content.Find.Execute("12×40=", $true, $false, $false, $false, $false, $true, 1, $false, "34×62=", 2) | Out-Null
$d.Content.Find.Execute("64×65=", $true, $false, $false, $false, $false, $true, 1, $false, "26×16=", 2) | Out-Null
$d.Content.Find.Execute("27×28=", $true, $false, $false, $false, $false, $true, 1, $false, "39×77=", 2) | Out-Null
